$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits after "33.79 PSU salinity".
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 2) Insert a new list paragraph "33.43 PSU" right after the paragraph that
#    reads "RMT25 nets – Pages 190-191" (JR177 section), before the
#    following hyperlink paragraph.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("RMT25 nets " + [char]8211 + " Pages 190-191")
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()

$newPara1 = $anchor.Paragraphs(1).Next()
$newPara1.Range.Text = "33.43 PSU"

# ---------------------------------------------------------------------------
# 3) Insert a new list paragraph "Use average salinity from other cruises
#    (33.65 PSU)." right after "Use midpoint 132" (JR38 section), before the
#    following hyperlink paragraph, and re-create the "_GoBack" bookmark as
#    a zero-length bookmark right at the end of that new paragraph's text.
# ---------------------------------------------------------------------------
$anchor2 = $d.Content
$anchor2.Find.Execute("Use midpoint 132")
$anchor2.Collapse(0)
$anchor2.InsertParagraphAfter()

$newPara2 = $anchor2.Paragraphs(1).Next()

$finalText = "Use average salinity from other cruises (33.65 PSU)."
# Type an extra sentinel character so the true end-of-text boundary is not
# the very last position in the paragraph (avoids a boundary case where a
# zero-length bookmark placed at the last position of a paragraph's final
# run gets mis-anchored). Add the bookmark at the boundary just before the
# sentinel, then delete the sentinel.
$newPara2.Range.Text = $finalText + "#"

$bmRange = $newPara2.Range.Duplicate
$bmRange.MoveEnd(1, -2)
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

$sentinel = $d.Range($newPara2.Range.End - 2, $newPara2.Range.End - 1)
$sentinel.Delete()
